$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '69.065.39'
$ws.Range("E2").Value = '  -3.90%  '

$ws.Range("D3").Value = '3.506.94'
$ws.Range("E3").Value = '  -4.99%  '

$ws.Range("E4").Value = '  -0.07%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '580.44'
$ws.Range("D5").Style = $ws.Range("C5").Style
$ws.Range("E5").Value = '  -1.59%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '174.30'
$ws.Range("D6").Style = $ws.Range("C6").Style
$ws.Range("E6").Value = '  -3.45%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.626'
$ws.Range("D7").Style = $ws.Range("C7").Style
$ws.Range("E7").Value = '  +0.63%  '

$ws.Range("D8").Value = '3.500.33'
$ws.Range("E8").Value = '  -4.95%  '

$ws.Range("E9").Value = '  +0.05%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.189'
$ws.Range("D10").Style = $ws.Range("C10").Style
$ws.Range("E10").Value = '  -6.21%  '

$ws.Range("E11").Value = '  +6.40%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.597'
$ws.Range("D12").Style = $ws.Range("C12").Style
$ws.Range("E12").Value = '  -3.10%  '

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '47.18'
$ws.Range("D13").Style = $ws.Range("C13").Style
$ws.Range("E13").Value = '  -5.87%  '

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '0.0000276'
$ws.Range("D14").Style = $ws.Range("C14").Style
$ws.Range("E14").Value = '  -3.77%  '

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '675.45'
$ws.Range("D15").Style = $ws.Range("C15").Style
$ws.Range("E15").Value = '  -1.75%  '

$ws.Range("D16").Value = '4.070.67'
$ws.Range("E16").Value = '  -5.11%  '

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '8.73'
$ws.Range("D17").Style = $ws.Range("C17").Style
$ws.Range("E17").Value = '  -3.36%  '

$ws.Range("B18").Value = 'WrappedBTC'
$ws.Range("C18").Value = 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc'
$ws.Range("D18").Value = '69.002.43'
$ws.Range("E18").Value = '  -4.07%  '

$ws.Range("B19").Value = 'WrappedEther'
$ws.Range("C19").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D19").Value = '3.509.83'
$ws.Range("E19").Value = '  -4.97%  '

$ws.Range("E20").Value = '  -1.37%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '17.53'
$ws.Range("D21").Style = $ws.Range("C21").Style
$ws.Range("E21").Value = '  -3.43%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '11.19'
$ws.Range("D22").Style = $ws.Range("C22").Style
$ws.Range("E22").Value = '  -4.39%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.905'
$ws.Range("D23").Style = $ws.Range("C23").Style
$ws.Range("E23").Value = '  -4.24%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '16.08'
$ws.Range("D24").Style = $ws.Range("C24").Style
$ws.Range("E24").Value = '  -10.03%  '

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '98.14'
$ws.Range("D25").Style = $ws.Range("C25").Style
$ws.Range("E25").Value = '  -5.81%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '3.86'
$ws.Range("D26").Style = $ws.Range("C26").Style
$ws.Range("E26").Value = '  -4.60%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '5.85'
$ws.Range("D27").Style = $ws.Range("C27").Style
$ws.Range("E27").Value = '  -0.06%  '

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '0.999'
$ws.Range("D28").Style = $ws.Range("C28").Style
$ws.Range("E28").Value = '  +0.01%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '2.66'
$ws.Range("D29").Style = $ws.Range("C29").Style
$ws.Range("E29").Value = '  -6.74%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '9.46'
$ws.Range("D30").Style = $ws.Range("C30").Style
$ws.Range("E30").Value = '  -7.17%  '

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '32.87'
$ws.Range("D31").Style = $ws.Range("C31").Style
$ws.Range("E31").Value = '  -7.65%  '

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '8.74'
$ws.Range("D32").Style = $ws.Range("C32").Style
$ws.Range("E32").Value = '  -5.47%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '3.21'
$ws.Range("D33").Style = $ws.Range("C33").Style
$ws.Range("E33").Value = '  -7.98%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '7.32'
$ws.Range("D34").Style = $ws.Range("C34").Style
$ws.Range("E34").Value = '  -0.76%  '

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.35'
$ws.Range("D35").Style = $ws.Range("C35").Style
$ws.Range("E35").Value = '  -5.72%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '590.31'
$ws.Range("D36").Style = $ws.Range("C36").Style
$ws.Range("E36").Value = '  +3.57%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '3.61'
$ws.Range("D37").Style = $ws.Range("C37").Style
$ws.Range("E37").Value = '  -15.77%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '10.91'
$ws.Range("D38").Style = $ws.Range("C38").Style
$ws.Range("E38").Value = '  -3.74%  '

$ws.Range("E39").Value = '  -4.42%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '57.27'
$ws.Range("D40").Style = $ws.Range("C40").Style
$ws.Range("E40").Value = '  -3.78%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.999'
$ws.Range("D41").Style = $ws.Range("C41").Style
$ws.Range("E41").Value = '  +0.09%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.0441'
$ws.Range("D42").Style = $ws.Range("C42").Style
$ws.Range("E42").Value = '  -5.78%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.338'
$ws.Range("D43").Style = $ws.Range("C43").Style
$ws.Range("E43").Value = '  -4.35%  '

$ws.Range("E44").Value = '  -6.95%  '

$ws.Range("D45").Value = '3.428.68'

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '33.46'
$ws.Range("D46").Style = $ws.Range("C46").Style
$ws.Range("E46").Value = '  -5.92%  '

$ws.Range("E47").Value = '  -9.58%  '

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '2.92'
$ws.Range("D48").Style = $ws.Range("C48").Style
$ws.Range("E48").Value = '  +0.80%  '

$ws.Range("E49").Value = '  -7.56%  '

$ws.Range("E50").Value = '  -0.40%  '

$ws.Range("B51").Value = 'MXToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '5.74'
$ws.Range("D51").Style = $ws.Range("C51").Style
$ws.Range("E51").Value = '  +17.63%  '
